$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DateText($cellRef, $text) {
    # Writing a "DD-MM-YYYY" string directly into a cell risks Excel's
    # autodetect reinterpreting it as a date serial number (e.g. when the
    # day portion is <= 12, it looks like a valid MM-DD-YYYY date). Route
    # the literal text through a scratch cell that's explicitly formatted
    # as Text, then copy/paste-values into the destination so the target
    # cell keeps its own formatting untouched and the text is never
    # reparsed as a date.
    $scratch = $ws.Range("AA1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $scratch.Clear()
}

# Row 3
Set-DateText "A3" "28-07-2022"
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# Row 4
Set-DateText "A4" "01-08-2022"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

# Row 5
Set-DateText "A5" "04-08-2022"

# Row 6
Set-DateText "A6" "08-08-2022"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("H6").Value = 0

# Row 7
Set-DateText "A7" "11-08-2022"

# Row 8
Set-DateText "A8" "15-08-2022"

# Row 9
Set-DateText "A9" "18-08-2022"

# Row 10
Set-DateText "A10" "22-08-2022"

# Row 11
Set-DateText "A11" "25-08-2022"
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 1
$ws.Range("H11").Value = 0

# Row 12
Set-DateText "A12" "29-08-2022"

# Row 13
Set-DateText "A13" "01-09-2022"
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("H13").Value = 0

# Row 14
Set-DateText "A14" "05-09-2022"

# Row 15
Set-DateText "A15" "08-09-2022"
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 1
$ws.Range("H15").Value = 0

# Row 16
Set-DateText "A16" "12-09-2022"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 1
$ws.Range("H16").Value = 0

# Row 17
Set-DateText "A17" "15-09-2022"

# Row 18
Set-DateText "A18" "19-09-2022"

# Row 19
Set-DateText "A19" "22-09-2022"

# Row 20
Set-DateText "A20" "26-09-2022"

# Row 21
Set-DateText "A21" "29-09-2022"
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 1
$ws.Range("H21").Value = 0
